$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "NA" values under the "duplicate_image_filename" column (column E) for rows 2-21
$ws.Range("E2:E21").Value = "NA"
